$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing contents/formatting so we can rebuild the table shifted
# one column to the left.
$ws.Cells.Clear()

# Header row (bold, centered, thin-bordered like the original header style)
$ws.Range("A1").Value = "QS_Astral_exact15"
$ws.Range("B1").Value = "FNRATE_ASTRAL"
$ws.Range("C1").Value = "TAXON"
$ws.Range("D1").Value = "MODEL_CONDITION"
$ws.Range("E1").Value = "GENE"

$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data rows
$ws.Range("A2").Value = 4148
$ws.Range("B2").Value = 0.125
$ws.Range("C2").Value = "11-texon"
$ws.Range("D2").Value = "estimated_15genes_strongILS"
$ws.Range("E2").Value = 11

$ws.Range("A3").Value = 4148
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "11-texon"
$ws.Range("D3").Value = "estimated_15genes_strongILS"
$ws.Range("E3").Value = 15
